$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet to reflect new "through" date
$ws.Name = "Through 2022-02-25"

# Update header label for the running total column
$ws.Range("I1").Value = "2022 (through 02-25)"

# Update February 2022 value (row 3, column I)
$ws.Range("I3").Value = 125

# Update the Total row (row 14, column I) to match new sum
$ws.Range("I14").Value = 284
